# MCQBert - full results.xlsx
# Fill in previously-blank accuracy / F1 / MCC cells across the three
# EPOCH tables on Sheet1, then leave the view scrolled/selected near the
# bottom of the sheet (K35), matching the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- EPOCH 1 table (rows 3-13): MISTRAL 10 (row7) & MISTRAL 40 (row11) ----
$ws.Range("H7").Value  = 0.783
$ws.Range("I7").Value  = 0.822
$ws.Range("J7").Value  = 0.722
$ws.Range("K7").Value  = 0.546

$ws.Range("H11").Value = 0.777
$ws.Range("I11").Value = 0.818
$ws.Range("J11").Value = 0.714
$ws.Range("K11").Value = 0.534

# ---- EPOCH 2 table (rows 17-27): MISTRAL 10 (row21) & MISTRAL 40 (row25) ----
$ws.Range("H21").Value = 0.791
$ws.Range("I21").Value = 0.824
$ws.Range("J21").Value = 0.742
$ws.Range("K21").Value = 0.566

$ws.Range("H25").Value = 0.79
$ws.Range("I25").Value = 0.83
$ws.Range("J25").Value = 0.727
$ws.Range("K25").Value = 0.561

# ---- EPOCH 3 table (rows 31-41): MISTRAL 10 (row35) & MISTRAL 40 (row39) ----
$ws.Range("C35").Value = 0.792
$ws.Range("D35").Value = 0.826
$ws.Range("E35").Value = 0.741
$ws.Range("F35").Value = 0.567
$ws.Range("H35").Value = 0.774
$ws.Range("I35").Value = 0.825
$ws.Range("J35").Value = 0.682
$ws.Range("K35").Value = 0.528

$ws.Range("C39").Value = 0.795
$ws.Range("D39").Value = 0.825
$ws.Range("E39").Value = 0.753
$ws.Range("F39").Value = 0.578
$ws.Range("H39").Value = 0.785
$ws.Range("I39").Value = 0.826
$ws.Range("J39").Value = 0.718
$ws.Range("K39").Value = 0.549

# ---- Leave the workbook scrolled to / selecting the last cell touched ----
$ws.Activate()
$ws.Range("K35").Select()
